$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit message indicates this is a plain SAVE of the project; the only
# functional change captured by the diff is the value stored in C10
# (rule "R30" / column "C1") changing from 18 to 100.
$ws.Range("C10").Value = 100
